$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Cells.Item(116, 8).Value = 31257304
$ws.Cells.Item(116, 9).Value = 62502810
$ws.Cells.Item(116, 10).Value = 11801.5
$ws.Cells.Item(116, 11).Value = 62502810
$ws.Cells.Item(116, 12).Value = 11801.5
$ws.Cells.Item(116, 13).Value = -62499368
$ws.Cells.Item(116, 14).Value = -18685.5
# Row 138
$ws.Cells.Item(138, 8).Value = 4238.029
$ws.Cells.Item(138, 9).Value = 1329.75
$ws.Cells.Item(138, 10).Value = 6687.1055
$ws.Cells.Item(138, 11).Value = 3989.25
$ws.Cells.Item(138, 12).Value = 20061.3165
$ws.Cells.Item(138, 13).Value = 1150.75
$ws.Cells.Item(138, 14).Value = -30341.3165

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Cells.Item(102, 8).Value = 1057.8
$ws.Cells.Item(102, 9).Value = 1057.8
$ws.Cells.Item(102, 11).Value = 1057.8
$ws.Cells.Item(102, 13).Value = 564.2
# Row 124
$ws.Cells.Item(124, 8).Value = 55799
$ws.Cells.Item(124, 10).Value = 55799
$ws.Cells.Item(124, 12).Value = 55799
$ws.Cells.Item(124, 14).Value = -65619

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Cells.Item(7, 8).Value = 1500
$ws.Cells.Item(7, 9).Value = 1500
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 1500
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -1387
$ws.Cells.Item(7, 14).Value = $null
# Row 16
$ws.Cells.Item(16, 8).Value = 1000
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 12).Value = 1000
$ws.Cells.Item(16, 14).Value = -1340
# Row 86
$ws.Cells.Item(86, 8).Value = 45504564
$ws.Cells.Item(86, 9).Value = 75220.07000000001
$ws.Cells.Item(86, 11).Value = 75220.07000000001
$ws.Cells.Item(86, 13).Value = -74097.07000000001
# Row 89
$ws.Cells.Item(89, 8).Value = 45504564
$ws.Cells.Item(89, 9).Value = 75220.07000000001
$ws.Cells.Item(89, 11).Value = 376100.35
$ws.Cells.Item(89, 13).Value = -370484.35
# Row 99
$ws.Cells.Item(99, 8).Value = 3790901
$ws.Cells.Item(99, 9).Value = 2943.9375
$ws.Cells.Item(99, 10).Value = 11366815
$ws.Cells.Item(99, 11).Value = 2943.9375
$ws.Cells.Item(99, 12).Value = 11366815
$ws.Cells.Item(99, 13).Value = -1445.9375
$ws.Cells.Item(99, 14).Value = -11369811

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 231.8
$ws.Cells.Item(22, 9).Value = 213.66667
$ws.Cells.Item(22, 10).Value = 259
$ws.Cells.Item(22, 11).Value = 213.66667
$ws.Cells.Item(22, 12).Value = 259
$ws.Cells.Item(22, 13).Value = 136.33333
$ws.Cells.Item(22, 14).Value = -959
# Row 29
$ws.Cells.Item(29, 8).Value = 2510.1667
$ws.Cells.Item(29, 10).Value = 3360.5
$ws.Cells.Item(29, 12).Value = 3360.5
$ws.Cells.Item(29, 14).Value = -3946.5
# Row 31
$ws.Cells.Item(31, 8).Value = 6218.582
$ws.Cells.Item(31, 9).Value = 2316.6
$ws.Cells.Item(31, 10).Value = 13047.05
$ws.Cells.Item(31, 11).Value = 2316.6
$ws.Cells.Item(31, 12).Value = 13047.05
$ws.Cells.Item(31, 13).Value = -2021.6
$ws.Cells.Item(31, 14).Value = -13637.05
# Row 34
$ws.Cells.Item(34, 8).Value = 6218.582
$ws.Cells.Item(34, 9).Value = 2316.6
$ws.Cells.Item(34, 10).Value = 13047.05
$ws.Cells.Item(34, 11).Value = 2316.6
$ws.Cells.Item(34, 12).Value = 13047.05
$ws.Cells.Item(34, 13).Value = -2114.6
$ws.Cells.Item(34, 14).Value = -13451.05
# Row 36
$ws.Cells.Item(36, 8).Value = 49016
$ws.Cells.Item(36, 10).Value = 49495
$ws.Cells.Item(36, 12).Value = 49495
$ws.Cells.Item(36, 14).Value = -50271
# Row 40
$ws.Cells.Item(40, 8).Value = 49016
$ws.Cells.Item(40, 10).Value = 49495
$ws.Cells.Item(40, 12).Value = 49495
$ws.Cells.Item(40, 14).Value = -49815
# Row 57
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).Value = $null
# Row 59
$ws.Cells.Item(59, 8).Value = 83333.336
$ws.Cells.Item(59, 10).Value = 100000
$ws.Cells.Item(59, 12).Value = 100000
$ws.Cells.Item(59, 14).Value = -102290
# Row 86
$ws.Cells.Item(86, 8).Value = 10504041
$ws.Cells.Item(86, 9).Value = 15631062
$ws.Cells.Item(86, 11).Value = 15631062
$ws.Cells.Item(86, 13).Value = -15629939
# Row 89
$ws.Cells.Item(89, 8).Value = 10504041
$ws.Cells.Item(89, 9).Value = 15631062
$ws.Cells.Item(89, 11).Value = 78155310
$ws.Cells.Item(89, 13).Value = -78149694
# Row 124
$ws.Cells.Item(124, 8).Value = 84998
$ws.Cells.Item(124, 10).Value = 84998
$ws.Cells.Item(124, 12).Value = 84998
$ws.Cells.Item(124, 14).Value = -89908
# Row 132
$ws.Cells.Item(132, 8).Value = 3698.9492
$ws.Cells.Item(132, 9).Value = 1791.8536
$ws.Cells.Item(132, 10).Value = 8042.8887
$ws.Cells.Item(132, 11).Value = 5375.560799999999
$ws.Cells.Item(132, 12).Value = 24128.6661
$ws.Cells.Item(132, 13).Value = -2845.560799999999
$ws.Cells.Item(132, 14).Value = -29188.6661

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 2355182.5
$ws.Cells.Item(5, 9).Value = 3636873
$ws.Cells.Item(5, 10).Value = 5416.6665
$ws.Cells.Item(5, 11).Value = 10910619
$ws.Cells.Item(5, 12).Value = 16249.9995
$ws.Cells.Item(5, 13).Value = -10910507
$ws.Cells.Item(5, 14).Value = -16473.9995
# Row 44
$ws.Cells.Item(44, 8).Value = 1201.25
$ws.Cells.Item(44, 9).Value = 286.66666
$ws.Cells.Item(44, 10).Value = 1750
$ws.Cells.Item(44, 11).Value = 859.9999799999999
$ws.Cells.Item(44, 12).Value = 5250
$ws.Cells.Item(44, 13).Value = -461.9999799999999
$ws.Cells.Item(44, 14).Value = -6046
# Row 56
$ws.Cells.Item(56, 8).Value = 7499
$ws.Cells.Item(56, 9).Value = 7499
$ws.Cells.Item(56, 11).Value = 7499
$ws.Cells.Item(56, 13).Value = -6969
# Row 81
$ws.Cells.Item(81, 8).Value = 6701
$ws.Cells.Item(81, 10).Value = 5999.8335
$ws.Cells.Item(81, 12).Value = 17999.5005
$ws.Cells.Item(81, 14).Value = -20245.5005
# Row 84
$ws.Cells.Item(84, 8).Value = 6701
$ws.Cells.Item(84, 10).Value = 5999.8335
$ws.Cells.Item(84, 12).Value = 53998.5015
$ws.Cells.Item(84, 14).Value = -65230.5015
# Row 124
$ws.Cells.Item(124, 8).Value = 3542.25
$ws.Cells.Item(124, 9).Value = 3542.25
$ws.Cells.Item(124, 11).Value = 10626.75
$ws.Cells.Item(124, 13).Value = -5716.75
# Row 126
$ws.Cells.Item(126, 8).Value = 1498.4
$ws.Cells.Item(126, 9).Value = 1139.75
$ws.Cells.Item(126, 11).Value = 3419.25
$ws.Cells.Item(126, 13).Value = 1520.75
# Row 129
$ws.Cells.Item(129, 8).Value = 1576.8572
$ws.Cells.Item(129, 9).Value = 1132
$ws.Cells.Item(129, 10).Value = 1824
$ws.Cells.Item(129, 11).Value = 3396
$ws.Cells.Item(129, 12).Value = 5472
$ws.Cells.Item(129, 13).Value = 1604
$ws.Cells.Item(129, 14).Value = -15472
# Row 131
$ws.Cells.Item(131, 8).Value = 1791.0625
$ws.Cells.Item(131, 10).Value = 1813.0834
$ws.Cells.Item(131, 12).Value = 5439.2502
$ws.Cells.Item(131, 14).Value = -15519.2502
# Row 135
$ws.Cells.Item(135, 8).Value = 2355182.5
$ws.Cells.Item(135, 9).Value = 3636873
$ws.Cells.Item(135, 10).Value = 5416.6665
$ws.Cells.Item(135, 11).Value = 32731857
$ws.Cells.Item(135, 12).Value = 48749.9985
$ws.Cells.Item(135, 13).Value = -32729322
$ws.Cells.Item(135, 14).Value = -53819.9985

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Cells.Item(20, 8).Value = 1733328.4
$ws.Cells.Item(20, 10).Value = 2500000
$ws.Cells.Item(20, 12).Value = 2500000
$ws.Cells.Item(20, 14).Value = -2500452
# Row 40
$ws.Cells.Item(40, 8).Value = 7274.091
$ws.Cells.Item(40, 9).Value = 4252.5
$ws.Cells.Item(40, 10).Value = 9000.714
$ws.Cells.Item(40, 11).Value = 4252.5
$ws.Cells.Item(40, 12).Value = 9000.714
$ws.Cells.Item(40, 13).Value = -4116.5
$ws.Cells.Item(40, 14).Value = -9272.714
# Row 46
$ws.Cells.Item(46, 8).Value = 2013.5714
$ws.Cells.Item(46, 9).Value = 669
$ws.Cells.Item(46, 10).Value = 2629.8333
$ws.Cells.Item(46, 11).Value = 669
$ws.Cells.Item(46, 12).Value = 2629.8333
$ws.Cells.Item(46, 13).Value = -481
$ws.Cells.Item(46, 14).Value = -3005.8333
# Row 132
$ws.Cells.Item(132, 8).Value = 7697178.5
$ws.Cells.Item(132, 9).Value = 13890596
$ws.Cells.Item(132, 11).Value = 41671788
$ws.Cells.Item(132, 13).Value = -41669258
